$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply cell updates per the cryptos list refresh (GitHub Actions data update).
# Numeric-looking Price values are prefixed with a literal apostrophe so Excel
# stores them as text (matching the source data, which keeps values like
# "97.10" / "0.0817" verbatim instead of coercing to a float and losing the
# trailing zero / fixed precision).

# Row 2
$ws.Cells.Item(2, 4).Value = "46.955.21"
$ws.Cells.Item(2, 5).Value = "  +5.64%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "2.336.35"
$ws.Cells.Item(3, 5).Value = "  +5.12%  "

# Row 4
$ws.Cells.Item(4, 5).Value = "  -0.87%  "

# Row 5
$ws.Cells.Item(5, 4).Value = "'305.83"
$ws.Cells.Item(5, 5).Value = "  +1.10%  "

# Row 6
$ws.Cells.Item(6, 4).Value = "'97.10"
$ws.Cells.Item(6, 5).Value = "  +7.47%  "

# Row 7
$ws.Cells.Item(7, 5).Value = "  +3.36%  "

# Row 8
$ws.Cells.Item(8, 5).Value = "  -0.72%  "

# Row 9
$ws.Cells.Item(9, 4).Value = "'0.539"
$ws.Cells.Item(9, 5).Value = "  +8.18%  "

# Row 10
$ws.Cells.Item(10, 4).Value = "'35.99"
$ws.Cells.Item(10, 5).Value = "  +6.58%  "

# Row 11
$ws.Cells.Item(11, 4).Value = "'0.0813"
$ws.Cells.Item(11, 5).Value = "  +4.02%  "

# Row 12
$ws.Cells.Item(12, 4).Value = "'7.46"
$ws.Cells.Item(12, 5).Value = "  +7.27%  "

# Row 13
$ws.Cells.Item(13, 5).Value = "  -0.06%  "

# Row 14
$ws.Cells.Item(14, 4).Value = "2.693.04"
$ws.Cells.Item(14, 5).Value = "  +4.94%  "

# Row 15
$ws.Cells.Item(15, 4).Value = "2.341.36"
$ws.Cells.Item(15, 5).Value = "  +5.16%  "

# Row 16
$ws.Cells.Item(16, 4).Value = "'14.18"
$ws.Cells.Item(16, 5).Value = "  +7.45%  "

# Row 17
$ws.Cells.Item(17, 4).Value = "'0.838"
$ws.Cells.Item(17, 5).Value = "  +3.72%  "

# Row 18
$ws.Cells.Item(18, 4).Value = "46.820.35"
$ws.Cells.Item(18, 5).Value = "  +5.32%  "

# Row 19
$ws.Cells.Item(19, 4).Value = "'13.74"
$ws.Cells.Item(19, 5).Value = "  +21.21%  "

# Row 20
$ws.Cells.Item(20, 4).Value = "0.0₃0952"
$ws.Cells.Item(20, 5).Value = "  +4.55%  "

# Row 21
$ws.Cells.Item(21, 4).Value = "'6.21"
$ws.Cells.Item(21, 5).Value = "  +2.95%  "

# Row 22
$ws.Cells.Item(22, 4).Value = "'67.80"
$ws.Cells.Item(22, 5).Value = "  +5.14%  "

# Row 23
$ws.Cells.Item(23, 4).Value = "'251.83"
$ws.Cells.Item(23, 5).Value = "  +7.85%  "

# Row 24
$ws.Cells.Item(24, 5).Value = "  +4.10%  "

# Row 25
$ws.Cells.Item(25, 5).Value = "  +4.72%  "

# Row 26
$ws.Cells.Item(26, 5).Value = "  -0.41%  "

# Row 27
$ws.Cells.Item(27, 4).Value = "'42.62"
$ws.Cells.Item(27, 5).Value = "  +17.70%  "

# Row 28
$ws.Cells.Item(28, 5).Value = "  -0.24%  "

# Row 29
$ws.Cells.Item(29, 5).Value = "  +4.63%  "

# Row 30
$ws.Cells.Item(30, 4).Value = "'20.36"
$ws.Cells.Item(30, 5).Value = "  +4.27%  "

# Row 31
$ws.Cells.Item(31, 4).Value = "'5.84"

# Row 32
$ws.Cells.Item(32, 4).Value = "'0.0817"
$ws.Cells.Item(32, 5).Value = "  +7.94%  "

# Row 33
$ws.Cells.Item(33, 4).Value = "'147.06"
$ws.Cells.Item(33, 5).Value = "  +0.59%  "

# Row 34
$ws.Cells.Item(34, 4).Value = "'2.62"

# Row 35
$ws.Cells.Item(35, 2).Value = "Kaspa"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Cells.Item(35, 4).Value = "'0.115"
$ws.Cells.Item(35, 5).Value = "  +8.39%  "

# Row 36
$ws.Cells.Item(36, 2).Value = "LidoDAOToken"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Cells.Item(36, 4).Value = "'3.15"
$ws.Cells.Item(36, 5).Value = "  +4.24%  "

# Row 37
$ws.Cells.Item(37, 4).Value = "'0.119"
$ws.Cells.Item(37, 5).Value = "  +3.09%  "

# Row 38
$ws.Cells.Item(38, 5).Value = "  +1.99%  "

# Row 39
$ws.Cells.Item(39, 5).Value = "  +10.46%  "

# Row 40
$ws.Cells.Item(40, 4).Value = "'0.0312"
$ws.Cells.Item(40, 5).Value = "  +8.18%  "

# Row 41
$ws.Cells.Item(41, 5).Value = "  +5.13%  "

# Row 42
$ws.Cells.Item(42, 4).Value = "'14.05"
$ws.Cells.Item(42, 5).Value = "  -3.10%  "

# Row 43
$ws.Cells.Item(43, 5).Value = "  -0.96%  "

# Row 44
$ws.Cells.Item(44, 4).Value = "'1.97"
$ws.Cells.Item(44, 5).Value = "  +14.77%  "

# Row 45
$ws.Cells.Item(45, 2).Value = "Maker"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Cells.Item(45, 4).Value = "1.811.86"
$ws.Cells.Item(45, 5).Value = "  +1.61%  "

# Row 46
$ws.Cells.Item(46, 2).Value = "BitcoinSV"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Cells.Item(46, 4).Value = "'91.19"
$ws.Cells.Item(46, 5).Value = "  +15.13%  "

# Row 47
$ws.Cells.Item(47, 4).Value = "'75.31"
$ws.Cells.Item(47, 5).Value = "  +12.47%  "

# Row 48
$ws.Cells.Item(48, 5).Value = "  +7.84%  "

# Row 49
$ws.Cells.Item(49, 4).Value = "'99.25"
$ws.Cells.Item(49, 5).Value = "  +4.00%  "

# Row 50
$ws.Cells.Item(50, 4).Value = "'55.47"
$ws.Cells.Item(50, 5).Value = "  +5.85%  "

# Row 51
$ws.Cells.Item(51, 2).Value = "FraxShare"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Cells.Item(51, 4).Value = "'8.08"
$ws.Cells.Item(51, 5).Value = "  +6.03%  "
